$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-9 gets updated from serial date 45233
# (2023-11-03) to serial date 45243 (2023-11-13).
for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = (Get-Date -Year 2023 -Month 11 -Day 13 -Hour 0 -Minute 0 -Second 0).Date
}
